$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 1111, shifting existing rows
# 1111-1209 down to 1114-1212 (dimension grows from A1:T1209 to A1:T1212).
$ws.Rows.Item(1111).Resize(3).Insert()

# Fill in the newly inserted rows with a new weekly price record
# (Terminal La Palmera de La Serena / Platano), matching the pattern used
# by every other triplet of rows in this table (Pinton / Primera Maduro /
# Primera Pinton for the same fecha).

# Row 1111: Pinton
$ws.Cells.Item(1111,1).Value = 8
$ws.Cells.Item(1111,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1111,3).Value = "Coquimbo"
$ws.Cells.Item(1111,4).Value = 45106
$ws.Cells.Item(1111,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1111,5).Value = 4
$ws.Cells.Item(1111,6).Value = "Fruta"
$ws.Cells.Item(1111,7).Value = 100108
$ws.Cells.Item(1111,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(1111,9).Value = 100108006
$ws.Cells.Item(1111,10).Value = "Plátano"
$ws.Cells.Item(1111,11).Value = "Sin especificar"
$ws.Cells.Item(1111,12).Value = "Pintón"
$ws.Cells.Item(1111,13).Value = 80
$ws.Cells.Item(1111,14).Value = 13000
$ws.Cells.Item(1111,15).Value = 13000
$ws.Cells.Item(1111,16).Value = 13000
$ws.Cells.Item(1111,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1111,18).Value = "Ecuador"
$ws.Cells.Item(1111,19).Value = 650
$ws.Cells.Item(1111,20).Value = 20

# Row 1112: Primera Maduro
$ws.Cells.Item(1112,1).Value = 8
$ws.Cells.Item(1112,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1112,3).Value = "Coquimbo"
$ws.Cells.Item(1112,4).Value = 45106
$ws.Cells.Item(1112,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1112,5).Value = 4
$ws.Cells.Item(1112,6).Value = "Fruta"
$ws.Cells.Item(1112,7).Value = 100108
$ws.Cells.Item(1112,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(1112,9).Value = 100108006
$ws.Cells.Item(1112,10).Value = "Plátano"
$ws.Cells.Item(1112,11).Value = "Sin especificar"
$ws.Cells.Item(1112,12).Value = "Primera Maduro"
$ws.Cells.Item(1112,13).Value = 120
$ws.Cells.Item(1112,14).Value = 15000
$ws.Cells.Item(1112,15).Value = 15000
$ws.Cells.Item(1112,16).Value = 15000
$ws.Cells.Item(1112,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1112,18).Value = "Ecuador"
$ws.Cells.Item(1112,19).Value = 750
$ws.Cells.Item(1112,20).Value = 20

# Row 1113: Primera Pintón
$ws.Cells.Item(1113,1).Value = 8
$ws.Cells.Item(1113,2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1113,3).Value = "Coquimbo"
$ws.Cells.Item(1113,4).Value = 45106
$ws.Cells.Item(1113,4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1113,5).Value = 4
$ws.Cells.Item(1113,6).Value = "Fruta"
$ws.Cells.Item(1113,7).Value = 100108
$ws.Cells.Item(1113,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(1113,9).Value = 100108006
$ws.Cells.Item(1113,10).Value = "Plátano"
$ws.Cells.Item(1113,11).Value = "Sin especificar"
$ws.Cells.Item(1113,12).Value = "Primera Pintón"
$ws.Cells.Item(1113,13).Value = 120
$ws.Cells.Item(1113,14).Value = 16000
$ws.Cells.Item(1113,15).Value = 16000
$ws.Cells.Item(1113,16).Value = 16000
$ws.Cells.Item(1113,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(1113,18).Value = "Ecuador"
$ws.Cells.Item(1113,19).Value = 800
$ws.Cells.Item(1113,20).Value = 20
